$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F325").Value = 775069
$ws.Range("F332").Value = 484997
$ws.Range("F337").Value = 105094
$ws.Range("F339").Value = 662447
$ws.Range("F349").Value = 159672
$ws.Range("G349").Value = 2758
$ws.Range("F352").Value = 306902
$ws.Range("F360").Value = 749294
$ws.Range("G360").Value = 5136
$ws.Range("F366").Value = 338784
$ws.Range("G366").Value = 2839
$ws.Range("F373").Value = 350144
$ws.Range("F376").Value = 222150
$ws.Range("F388").Value = 730441
$ws.Range("G388").Value = 2202
$ws.Range("F421").Value = 153287
$ws.Range("F422").Value = 298600
$ws.Range("G422").Value = 646
$ws.Range("F423").Value = 438732
$ws.Range("F425").Value = 138368
$ws.Range("F426").Value = 106951
$ws.Range("F431").Value = 170843
$ws.Range("F432").Value = 124045
$ws.Range("F433").Value = 86464
$ws.Range("G433").Value = 265
$ws.Range("F434").Value = 79406
$ws.Range("F435").Value = 83054
$ws.Range("F436").Value = 144575
$ws.Range("F437").Value = 167284
$ws.Range("F438").Value = 121527
$ws.Range("F439").Value = 89252
$ws.Range("F440").Value = 73555
$ws.Range("F441").Value = 68185
$ws.Range("F443").Value = 106950
$ws.Range("F444").Value = 103886
$ws.Range("G444").Value = 191
$ws.Range("F445").Value = 84585
$ws.Range("F446").Value = 86083
$ws.Range("F447").Value = 67095
$ws.Range("F448").Value = 61309
$ws.Range("G448").Value = 138
$ws.Range("F449").Value = 59713
$ws.Range("F450").Value = 91059
$ws.Range("F451").Value = 85337
$ws.Range("F452").Value = 74572
$ws.Range("G452").Value = 126
$ws.Range("F453").Value = 69743
$ws.Range("G453").Value = 210
$ws.Range("F454").Value = 50958
$ws.Range("G454").Value = 127
